# ECOSCOM-4453 - Create group action to export records to excel
# Changes according to the pull request comments:
# the template's sample/placeholder value in A2 is removed, leaving the
# cell blank (its style is preserved) while the "title" header in A1 stays.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
